# Extend the yearly statistics table with two new columns: S (2021) and T
# (2022), matching the existing number/border/font formatting that column R
# already uses for each row (header row 4, data rows 5-14).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column R's formatting (font/border/number format) into the new
# columns S and T, for every row of the table (this also seeds placeholder
# values that are immediately overwritten below).
$ws.Range("R4:R14").Copy($ws.Range("S4:S14"))
$ws.Range("R4:R14").Copy($ws.Range("T4:T14"))

# Row 4: header years.
$ws.Range("S4").Value = 2021
$ws.Range("T4").Value = 2022

# Row 5.
$ws.Range("S5").Value = 2.5
$ws.Range("T5").Value = 2.6

# Row 6.
$ws.Range("S6").Value = 2.5
$ws.Range("T6").Value = 1.8

# Row 7.
$ws.Range("S7").Value = 1.6
$ws.Range("T7").Value = 2.6

# Row 8.
$ws.Range("S8").Value = 3.6
$ws.Range("T8").Value = 1.9

# Row 9.
$ws.Range("S9").Value = 5.8
$ws.Range("T9").Value = 3.9

# Row 10.
$ws.Range("S10").Value = 1.1000000000000001
$ws.Range("T10").Value = 3.2

# Row 11.
$ws.Range("S11").Value = 1.1000000000000001
$ws.Range("T11").Value = 3.3

# Row 12.
$ws.Range("S12").Value = 5.0999999999999996
$ws.Range("T12").Value = 2.5

# Row 13.
$ws.Range("S13").Value = 2.2999999999999998
$ws.Range("T13").Value = 1.9

# Row 14.
$ws.Range("S14").Value = 2.1
$ws.Range("T14").Value = 2.5

# Match the author's final cell selection.
[void]$ws.Range("V7").Select()
